$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 111782493
$ws.Range("B2").Value = 98535
$ws.Range("E2").Value = 222498
$ws.Range("F2").Value = "Blåsippa"
$ws.Range("G2").Value = "Hepatica nobilis"
$ws.Range("H2").Value = "Schreb."
$ws.Range("I2").Value = "'10"
$ws.Range("Q2").Value = 574992
$ws.Range("R2").Value = 6299306
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# Row 3 updates
$ws.Range("A3").Value = 111782491
$ws.Range("B3").Value = 108022
$ws.Range("E3").Value = 219677
$ws.Range("F3").Value = "Murgröna"
$ws.Range("G3").Value = "Hedera helix"
$ws.Range("H3").Value = "L."
$ws.Range("I3").Value = "'1"
$ws.Range("Q3").Value = 575009
$ws.Range("R3").Value = 6299346
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
